$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.959.89"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").Value = "2.249.91"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.86"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.634"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "62.99"
$ws.Range("E7").Value = "  -1.76%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.447"
$ws.Range("E9").Value = "  +5.46%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0977"
$ws.Range("E10").Value = "  +2.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "57.18"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "26.30"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.105"
$ws.Range("E13").Value = "  +1.46%  "
$ws.Range("D14").Value = "2.582.23"
$ws.Range("E14").Value = "  -1.57%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.45"
$ws.Range("E15").Value = "  -1.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.04"
$ws.Range("E16").Value = "  +2.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.825"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "2.245.13"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "43.846.81"
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("E20").Value = "  +3.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.63"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "247.90"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("E25").Value = "  -6.04%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.29"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.32"
$ws.Range("E27").Value = "  +20.42%  "
$ws.Range("E28").Value = "  -0.42%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "173.11"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("E30").Value = "  +2.20%  "
$ws.Range("E31").Value = "  +0.26%  "
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0683"
$ws.Range("E34").Value = "  -2.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.80"
$ws.Range("E35").Value = "  +1.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.93"
$ws.Range("E36").Value = "  -3.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.65"
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.39"
$ws.Range("E38").Value = "  -3.41%  "
$ws.Range("E39").Value = "  -4.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0254"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.64"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.000224"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "17.04"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.48"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  -1.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0945"
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("E48").Value = "  -8.79%  "
$ws.Range("D49").Value = "1.438.54"
$ws.Range("E49").Value = "  -3.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.29"
$ws.Range("E50").Value = "  -2.55%  "
$ws.Range("E51").Value = "  +1.39%  "
